$d = $word.ActiveDocument

# The edit touches only the document's first paragraph (the hidden
# "**ID__...__ID**" bookmark/placeholder line right before the
# "SMC PGI 5345 Government Property" heading):
#   1. give the paragraph a border (just the 5-twip border/text spacing
#      on all four sides -- no visible line),
#   2. bump its left indent from 120 -> 225 twips,
#   3. rename the placeholder id from
#      **ID__AFFARS_pgi_5345_topic_8__ID** to
#      **ID__AFFARS_SMC_PGI_5345__ID**, and
#   4. drop the trailing run that used to hold a single, now
#      unnecessary, trailing space.

$p1 = $d.Paragraphs(1)

# --- 1. Paragraph border (renders as <w:pBdr><w:top w:space="5"/>...) ---
$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# --- 2. Left indent: 225 twips == 11.25 points (Word's ParagraphFormat
#        properties are expressed in points; OOXML w:ind is in twips). ---
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# --- 3 & 4. Swap the placeholder id text and remove the leftover
#           trailing-space run, locating everything via Find so we don't
#           depend on hard-coded character offsets. ---
$oldId = "**ID__AFFARS_pgi_5345_topic_8__ID**"
$newId = "**ID__AFFARS_SMC_PGI_5345__ID**"

$searchRange = $d.Range($p1.Range.Start, $p1.Range.End)
$found = $searchRange.Find.Execute($oldId, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $idRange = $d.Range($searchRange.Start, $searchRange.End)
    $oldIdEnd = $idRange.End

    # Replace just the placeholder run's text; this leaves its original
    # run formatting (rPr) untouched.
    $idRange.Text = $newId

    # Whatever used to sit immediately after the placeholder (here, a run
    # containing a single space) now starts `shift` characters earlier,
    # where `shift` is how many characters shorter the new id is.
    $shift = $oldId.Length - $newId.Length
    $trailingStart = $oldIdEnd - $shift

    # End just before the paragraph mark so we never touch it.
    $trailingEnd = $d.Paragraphs(1).Range.End - 1

    if ($trailingEnd -gt $trailingStart) {
        $trailingRange = $d.Range($trailingStart, $trailingEnd)
        if ($trailingRange.Text -eq " ") {
            $trailingRange.Text = ""
        }
    }
}
